$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 232.04347
$ws.Range("I33").Value = 225.4762
$ws.Range("K33").Value = 225.4762
$ws.Range("M33").Value = 3.523799999999994
$ws.Range("H40").Value = 1442.3889
$ws.Range("I40").Value = 1384.3334
$ws.Range("J40").Value = 1500.4445
$ws.Range("K40").Value = 1384.3334
$ws.Range("L40").Value = 1500.4445
$ws.Range("M40").Value = -1209.3334
$ws.Range("N40").Value = -1850.4445
$ws.Range("H64").Value = 3451.1738
$ws.Range("I64").Value = 3330.7778
$ws.Range("J64").Value = 3528.5715
$ws.Range("K64").Value = 3330.7778
$ws.Range("L64").Value = 3528.5715
$ws.Range("M64").Value = -3082.7778
$ws.Range("N64").Value = -4024.5715
$ws.Range("H67").Value = 3451.1738
$ws.Range("I67").Value = 3330.7778
$ws.Range("J67").Value = 3528.5715
$ws.Range("K67").Value = 3330.7778
$ws.Range("L67").Value = 3528.5715
$ws.Range("M67").Value = -2472.7778
$ws.Range("N67").Value = -5244.5715
$ws.Range("H70").Value = 8384596
$ws.Range("I70").Value = 27945360
$ws.Range("J70").Value = 1411.1428
$ws.Range("K70").Value = 83836080
$ws.Range("L70").Value = 4233.428400000001
$ws.Range("M70").Value = -83835810
$ws.Range("N70").Value = -4773.428400000001
$ws.Range("H73").Value = 8384596
$ws.Range("I73").Value = 27945360
$ws.Range("J73").Value = 1411.1428
$ws.Range("K73").Value = 83836080
$ws.Range("L73").Value = 4233.428400000001
$ws.Range("M73").Value = -83835144
$ws.Range("N73").Value = -6105.428400000001
$ws.Range("H74").Value = 3764.2942
$ws.Range("I74").Value = 3249.6667
$ws.Range("J74").Value = 4045
$ws.Range("K74").Value = 3249.6667
$ws.Range("L74").Value = 4045
$ws.Range("M74").Value = -2313.6667
$ws.Range("N74").Value = -5917
$ws.Range("H76").Value = 4485.8335
$ws.Range("I76").Value = 4590
$ws.Range("J76").Value = 4340
$ws.Range("K76").Value = 4590
$ws.Range("L76").Value = 4340
$ws.Range("M76").Value = -4275
$ws.Range("N76").Value = -4970
$ws.Range("H77").Value = 3764.2942
$ws.Range("I77").Value = 3249.6667
$ws.Range("J77").Value = 4045
$ws.Range("K77").Value = 16248.3335
$ws.Range("L77").Value = 20225
$ws.Range("M77").Value = -11568.3335
$ws.Range("N77").Value = -29585
$ws.Range("H79").Value = 4485.8335
$ws.Range("I79").Value = 4590
$ws.Range("J79").Value = 4340
$ws.Range("K79").Value = 4590
$ws.Range("L79").Value = 4340
$ws.Range("M79").Value = -3498
$ws.Range("N79").Value = -6524
$ws.Range("H81").Value = 35998.332
$ws.Range("J81").Value = 35998.332
$ws.Range("L81").Value = 35998.332
$ws.Range("N81").Value = -37994.332
$ws.Range("H84").Value = 35998.332
$ws.Range("J84").Value = 35998.332
$ws.Range("L84").Value = 107994.996
$ws.Range("N84").Value = -117978.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2638
$ws.Range("I45").Value = 2500
$ws.Range("J45").Value = 2914
$ws.Range("K45").Value = 2500
$ws.Range("L45").Value = 2914
$ws.Range("M45").Value = -2123
$ws.Range("N45").Value = -3668
$ws.Range("H88").Value = 2801.5
$ws.Range("I88").Value = 2068.6667
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 2068.6667
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -1662.6667
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 2801.5
$ws.Range("I91").Value = 2068.6667
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 2068.6667
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -664.6667000000002
$ws.Range("N91").Value = -7808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 78769.62
$ws.Range("I86").Value = 1999.8334
$ws.Range("J86").Value = 1000007
$ws.Range("K86").Value = 1999.8334
$ws.Range("L86").Value = 1000007
$ws.Range("M86").Value = -876.8334
$ws.Range("N86").Value = -1002253
$ws.Range("H89").Value = 78769.62
$ws.Range("I89").Value = 1999.8334
$ws.Range("J89").Value = 1000007
$ws.Range("K89").Value = 9999.166999999999
$ws.Range("L89").Value = 5000035
$ws.Range("M89").Value = -4383.166999999999
$ws.Range("N89").Value = -5011267
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 257.82608
$ws.Range("I22").Value = 215.78947
$ws.Range("J22").Value = 457.5
$ws.Range("K22").Value = 215.78947
$ws.Range("L22").Value = 457.5
$ws.Range("M22").Value = 134.21053
$ws.Range("N22").Value = -1157.5
$ws.Range("H62").Value = 74153.42999999999
$ws.Range("I62").Value = 86045.664
$ws.Range("J62").Value = 2800
$ws.Range("K62").Value = 86045.664
$ws.Range("L62").Value = 2800
$ws.Range("M62").Value = -85421.664
$ws.Range("N62").Value = -4048
$ws.Range("H65").Value = 74153.42999999999
$ws.Range("I65").Value = 86045.664
$ws.Range("J65").Value = 2800
$ws.Range("K65").Value = 430228.32
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = -427108.32
$ws.Range("N65").Value = -20240
$ws.Range("H74").Value = 29989.5
$ws.Range("J74").Value = 29989.5
$ws.Range("L74").Value = 29989.5
$ws.Range("N74").Value = -31737.5
$ws.Range("H77").Value = 29989.5
$ws.Range("J77").Value = 29989.5
$ws.Range("L77").Value = 89968.5
$ws.Range("N77").Value = -98704.5
$ws.Range("H107").Value = 377.2353
$ws.Range("I107").Value = 327.47827
$ws.Range("J107").Value = 481.27274
$ws.Range("K107").Value = 327.47827
$ws.Range("L107").Value = 481.27274
$ws.Range("M107").Value = 1592.52173
$ws.Range("N107").Value = -4321.27274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2125
$ws.Range("J54").Value = 2533.3333
$ws.Range("L54").Value = 7599.999899999999
$ws.Range("N54").Value = -8717.999899999999
$ws.Range("H139").Value = 1716.1072
$ws.Range("J139").Value = 3000
$ws.Range("L139").Value = 9000
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5852.5
$ws.Range("I70").Value = 5396.8
$ws.Range("J70").Value = 6829
$ws.Range("K70").Value = 5396.8
$ws.Range("L70").Value = 6829
$ws.Range("M70").Value = -5126.8
$ws.Range("N70").Value = -7369
$ws.Range("H73").Value = 5852.5
$ws.Range("I73").Value = 5396.8
$ws.Range("J73").Value = 6829
$ws.Range("K73").Value = 5396.8
$ws.Range("L73").Value = 6829
$ws.Range("M73").Value = -4460.8
$ws.Range("N73").Value = -8701
$ws.Range("H80").Value = 2678.2942
$ws.Range("I80").Value = 2535
$ws.Range("J80").Value = 2839.5
$ws.Range("K80").Value = 2535
$ws.Range("L80").Value = 2839.5
$ws.Range("M80").Value = -1537
$ws.Range("N80").Value = -4835.5
$ws.Range("H83").Value = 2678.2942
$ws.Range("I83").Value = 2535
$ws.Range("J83").Value = 2839.5
$ws.Range("K83").Value = 12675
$ws.Range("L83").Value = 14197.5
$ws.Range("M83").Value = -7683
$ws.Range("N83").Value = -24181.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1841.25
$ws.Range("I68").Value = 1455
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1455
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -706
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 1841.25
$ws.Range("I71").Value = 1455
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 7275
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -3531
$ws.Range("N71").Value = -22488
$ws.Range("H81").Value = 31181
$ws.Range("J81").Value = 31181
$ws.Range("L81").Value = 31181
$ws.Range("N81").Value = -33177
$ws.Range("H84").Value = 31181
$ws.Range("J84").Value = 31181
$ws.Range("L84").Value = 93543
$ws.Range("N84").Value = -103527

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 30000
$ws.Range("K76").Value = 30000
$ws.Range("M76").Value = -29685
$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 30000
$ws.Range("K79").Value = 30000
$ws.Range("M79").Value = -28908
